$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 21 new rows (2024-08-28 .. 2024-09-25) to the stock price history,
# mirroring the same date range that already appears earlier in the sheet
# but with freshly-computed returns (I) and NAV (J) chained off the prior
# last row (J638).

$ws.Range("A639").Value = "'2024-08-28"
$ws.Range("C639").Value = 1939.099975585938
$ws.Range("D639").Value = 1719.449951171875
$ws.Range("E639").Value = 1078.800048828125
$ws.Range("F639").Value = 1847.050048828125
$ws.Range("G639").Value = 1737.900024414062
$ws.Range("H639").Value = 37879.2001953125
$ws.Range("I639").Value = 0
$ws.Range("J639").Value = 173.8459610261267
$ws.Range("A640").Value = "'2024-08-29"
$ws.Range("C640").Value = 1933.349975585938
$ws.Range("D640").Value = 1751.849975585938
$ws.Range("E640").Value = 1061.300048828125
$ws.Range("F640").Value = 1857.849975585938
$ws.Range("G640").Value = 1726.550048828125
$ws.Range("H640").Value = 37876.40014648438
$ws.Range("I640").Value = -0.0000739204844264769429059
$ws.Range("J640").Value = 173.8331102484721
$ws.Range("A641").Value = "'2024-08-30"
$ws.Range("C641").Value = 1943.699951171875
$ws.Range("D641").Value = 1753.25
$ws.Range("E641").Value = 1065.599975585938
$ws.Range("F641").Value = 1815.150024414062
$ws.Range("G641").Value = 1772.25
$ws.Range("H641").Value = 37943.19958496094
$ws.Range("I641").Value = 0.001763616347335551
$ws.Range("J641").Value = 174.1396851634145
$ws.Range("A642").Value = "'2024-09-02"
$ws.Range("C642").Value = 1964.5
$ws.Range("D642").Value = 1806.650024414062
$ws.Range("E642").Value = 1050.949951171875
$ws.Range("F642").Value = 1766.300048828125
$ws.Range("G642").Value = 1749.5
$ws.Range("H642").Value = 38025.74987792969
$ws.Range("I642").Value = 0.002175628146063607
$ws.Range("J642").Value = 174.5185483638027
$ws.Range("A643").Value = "'2024-09-03"
$ws.Range("C643").Value = 1941.25
$ws.Range("D643").Value = 1790.449951171875
$ws.Range("E643").Value = 1068.800048828125
$ws.Range("F643").Value = 1769.300048828125
$ws.Range("G643").Value = 1718.75
$ws.Range("H643").Value = 37895.55029296875
$ws.Range("I643").Value = -0.003423984678248408
$ws.Range("J643").Value = 173.9209995281348
$ws.Range("A644").Value = "'2024-09-04"
$ws.Range("C644").Value = 1922.449951171875
$ws.Range("D644").Value = 1785.25
$ws.Range("E644").Value = 1056.199951171875
$ws.Range("F644").Value = 1749.699951171875
$ws.Range("G644").Value = 1729.550048828125
$ws.Range("H644").Value = 37618.74926757812
$ws.Range("I644").Value = -0.007304314708473397
$ws.Range("J644").Value = 172.6506258131691
$ws.Range("A645").Value = "'2024-09-05"
$ws.Range("C645").Value = 1933.150024414062
$ws.Range("D645").Value = 1790.550048828125
$ws.Range("E645").Value = 1074.900024414062
$ws.Range("F645").Value = 1722.900024414062
$ws.Range("G645").Value = 1720.75
$ws.Range("H645").Value = 37761.05065917969
$ws.Range("I645").Value = 0.003782725220059497
$ws.Range("J645").Value = 173.3037156896916
$ws.Range("A646").Value = "'2024-09-06"
$ws.Range("C646").Value = 1901.849975585938
$ws.Range("D646").Value = 1756.099975585938
$ws.Range("E646").Value = 1112.650024414062
$ws.Range("F646").Value = 1730.300048828125
$ws.Range("G646").Value = 1715
$ws.Range("H646").Value = 37713.70007324219
$ws.Range("I646").Value = -0.001253953084220899
$ws.Range("J646").Value = 173.0864009608956
$ws.Range("A647").Value = "'2024-09-09"
$ws.Range("C647").Value = 1894.650024414062
$ws.Range("D647").Value = 1746.75
$ws.Range("E647").Value = 1077.550048828125
$ws.Range("F647").Value = 1750.400024414062
$ws.Range("G647").Value = 1741.199951171875
$ws.Range("H647").Value = 37455.65051269531
$ws.Range("I647").Value = -0.006842329446480399
$ws.Range("J647").Value = 171.9020867828155
$ws.Range("A648").Value = "'2024-09-10"
$ws.Range("C648").Value = 1912.300048828125
$ws.Range("D648").Value = 1779.099975585938
$ws.Range("E648").Value = 1091
$ws.Range("F648").Value = 1756.349975585938
$ws.Range("G648").Value = 1745.150024414062
$ws.Range("H648").Value = 37856.65014648438
$ws.Range("I648").Value = 0.01070598503296977
$ws.Range("J648").Value = 173.7424679510486
$ws.Range("A649").Value = "'2024-09-11"
$ws.Range("C649").Value = 1910.150024414062
$ws.Range("D649").Value = 1778.75
$ws.Range("E649").Value = 1077.849975585938
$ws.Range("F649").Value = 1789.349975585938
$ws.Range("G649").Value = 1782.650024414062
$ws.Range("H649").Value = 37910.79992675781
$ws.Range("I649").Value = 0.001430390170918655
$ws.Range("J649").Value = 173.990987469477
$ws.Range("A650").Value = "'2024-09-12"
$ws.Range("C650").Value = 1950.449951171875
$ws.Range("D650").Value = 1807.599975585938
$ws.Range("E650").Value = 1083.75
$ws.Range("F650").Value = 1838.050048828125
$ws.Range("G650").Value = 1812.75
$ws.Range("H650").Value = 38550.34973144531
$ws.Range("I650").Value = 0.01686985782212682
$ws.Range("J650").Value = 176.9261906904185
$ws.Range("A651").Value = "'2024-09-13"
$ws.Range("C651").Value = 1944.099975585938
$ws.Range("D651").Value = 1812.800048828125
$ws.Range("E651").Value = 1089.699951171875
$ws.Range("F651").Value = 1826.050048828125
$ws.Range("G651").Value = 1814.099975585938
$ws.Range("H651").Value = 38552.5498046875
$ws.Range("I651").Value = 0.0000570701240718683331862
$ws.Range("J651").Value = 176.9362878900728
$ws.Range("A652").Value = "'2024-09-16"
$ws.Range("C652").Value = 1950.25
$ws.Range("D652").Value = 1811.849975585938
$ws.Range("E652").Value = 1094.650024414062
$ws.Range("F652").Value = 1757.849975585938
$ws.Range("G652").Value = 1797.199951171875
$ws.Range("H652").Value = 38385.89990234375
$ws.Range("I652").Value = -0.004322668751821118
$ws.Range("J652").Value = 176.1714509273471
$ws.Range("A653").Value = "'2024-09-17"
$ws.Range("C653").Value = 1952.550048828125
$ws.Range("D653").Value = 1813.75
$ws.Range("E653").Value = 1080.300048828125
$ws.Range("F653").Value = 1741.150024414062
$ws.Range("G653").Value = 1848.5
$ws.Range("H653").Value = 38346.90075683594
$ws.Range("I653").Value = -0.001015975803798501
$ws.Range("J653").Value = 175.9924649958849
$ws.Range("A654").Value = "'2024-09-18"
$ws.Range("C654").Value = 1892.150024414062
$ws.Range("D654").Value = 1756.5
$ws.Range("E654").Value = 1065.800048828125
$ws.Range("F654").Value = 1727.25
$ws.Range("G654").Value = 1805.599975585938
$ws.Range("H654").Value = 37454.75048828125
$ws.Range("I654").Value = -0.023265250931541
$ws.Range("J654").Value = 171.8979561356952
$ws.Range("A655").Value = "'2024-09-19"
$ws.Range("C655").Value = 1894.199951171875
$ws.Range("D655").Value = 1736.5
$ws.Range("E655").Value = 1060.75
$ws.Range("F655").Value = 1676.449951171875
$ws.Range("G655").Value = 1877.449951171875
$ws.Range("H655").Value = 37317.94946289062
$ws.Range("I655").Value = -0.003652434567236724
$ws.Range("J655").Value = 171.2701100986678
$ws.Range("A656").Value = "'2024-09-20"
$ws.Range("C656").Value = 1905.75
$ws.Range("D656").Value = 1760.050048828125
$ws.Range("E656").Value = 1114.699951171875
$ws.Range("F656").Value = 1662
$ws.Range("G656").Value = 1931.449951171875
$ws.Range("H656").Value = 38001.24975585938
$ws.Range("I656").Value = 0.01831023147850691
$ws.Range("J656").Value = 174.4061054599238
$ws.Range("A657").Value = "'2024-09-23"
$ws.Range("C657").Value = 1896.449951171875
$ws.Range("D657").Value = 1752.800048828125
$ws.Range("E657").Value = 1106.699951171875
$ws.Range("F657").Value = 1692.900024414062
$ws.Range("G657").Value = 1920.400024414062
$ws.Range("H657").Value = 37915.79968261719
$ws.Range("I657").Value = -0.002248612184892999
$ws.Range("J657").Value = 174.0139337660669
$ws.Range("A658").Value = "'2024-09-24"
$ws.Range("C658").Value = 1898.599975585938
$ws.Range("D658").Value = 1775.599975585938
$ws.Range("E658").Value = 1098.5
$ws.Range("F658").Value = 1660.900024414062
$ws.Range("G658").Value = 1838.75
$ws.Range("H658").Value = 37717.7998046875
$ws.Range("I658").Value = -0.00522209420840627
$ws.Range("J658").Value = 173.1052166103651
$ws.Range("A659").Value = "'2024-09-25"
$ws.Range("C659").Value = 1895.300048828125
$ws.Range("D659").Value = 1782.400024414062
$ws.Range("E659").Value = 1088.599975585938
$ws.Range("F659").Value = 1654.75
$ws.Range("G659").Value = 1722.050048828125
$ws.Range("H659").Value = 37400.95031738281
$ws.Range("I659").Value = -0.008400529430280024
$ws.Range("J659").Value = 171.6510411436947

# The literal date strings in column A were entered with a leading
# apostrophe above so Excel keeps them as text (matching the rest of
# column A) instead of auto-converting to date serials. That also stamps
# a "quote prefix" style on those cells, so normalise their formatting
# back to the same (default) style already used by the rest of column A.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A639:A659").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
